# Add a "(After 30 secs deadband)" note under the "VSENSE_IN = H" label that
# sits in the TV_OFF part of the power-control state diagram, and resize the
# textbox so the new second line fits without wrapping.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# There are two "VSENSE_IN = H" textboxes on the slide; the one we want is
# the lower one, positioned at roughly (4654786, 2821421) EMU
# (~366.52pt, ~222.16pt) next to the TV_OFF state box.
$targetLeftPt = 4654786 / 12700
$targetTopPt = 2821421 / 12700

$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText -and `
        $candidate.TextFrame.TextRange.Text -eq "VSENSE_IN = H" -and `
        [Math]::Abs($candidate.Left - $targetLeftPt) -lt 1 -and `
        [Math]::Abs($candidate.Top - $targetTopPt) -lt 1) {
        $shape = $candidate
        break
    }
}

if ($shape -eq $null) {
    throw "Could not find the VSENSE_IN = H textbox near TV_OFF"
}

# Append a second paragraph with the deadband note.
$shape.TextFrame.TextRange.Text = "VSENSE_IN = H`r(After 30 secs deadband)"

# Resize/reposition the textbox to its new (autofit) bounding box, matching
# the wider text now that a second line has been added.
$shape.Left = 365.5508661417323
$shape.Top = 222.15913385826772
$shape.Width = 157.68070866141733
$shape.Height = 33.92811023622047
